# [Fonds de solidarite] Add 2020-07-21 data
# Updates "nombre_aides" (col C) and "montant_total" (col D) for the rows
# affected by the 2020-07-21 refresh of the regional / classe-effectif
# breakdown. Values are written back as text (matching the sheet's existing
# inline-string storage for these numeric-looking columns), so the number
# format is forced to "@" (Text) before assignment to avoid Excel
# auto-coercing the strings to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (nombre_aides, montant_total)
$updates = @{
    9  = @("36",   "82250.00")
    10 = @("262",  "631590.66")
    11 = @("109",  "334891.77")
    15 = @("85",   "204152.38")
    16 = @("373",  "1041917.16")
    17 = @("125",  "453500.00")
    18 = @("38",   "156045.00")
    19 = @("10",   "63850.00")
    44 = @("9",    "27500.00")
    45 = @("35",   "137903.07")
    48 = @("69",   "192718.00")
    49 = @("458",  "1309380.59")
    50 = @("200",  "695300.15")
    51 = @("58",   "279851.23")
    52 = @("18",   "90000.00")
    53 = @("12",   "30720.65")
    55 = @("2482", "5818748.44")
    56 = @("1281", "3721348.00")
    71 = @("197",  "478326.09")
    72 = @("771",  "2135863.23")
    73 = @("288",  "993266.79")
    74 = @("92",   "364100.00")
    75 = @("21",   "115383.20")
    76 = @("21",   "43500.00")
    77 = @("181",  "403500.00")
    78 = @("801",  "1892198.76")
    81 = @("13",   "53000.00")
    82 = @("24",   "48000.00")
}

foreach ($row in $updates.Keys) {
    $values = $updates[$row]

    $cellC = $ws.Range("C$row")
    $cellC.NumberFormat = "@"
    $cellC.Value = $values[0]

    $cellD = $ws.Range("D$row")
    $cellD.NumberFormat = "@"
    $cellD.Value = $values[1]
}
